# Insert a new weekly price record for Coliflor at Vega Monumental Concepción.
# This pushes the existing rows 106..205 down to 107..206 (the former last
# row, 205, becomes row 206), and the freshly inserted row 106 receives the
# new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 106..205 down by one to make room for the new record.
$ws.Rows(106).Insert()

# Populate the newly inserted row 106 with the new record's values.
$ws.Range("A106").Value2 = 11
$ws.Range("B106").Value2 = "Vega Monumental Concepción"
$ws.Range("C106").Value2 = "Bíobío"
$ws.Range("D106").Value2 = 44586
$ws.Range("E106").Value2 = 8
$ws.Range("F106").Value2 = 100112008
$ws.Range("G106").Value2 = "Coliflor"
$ws.Range("H106").Value2 = "Sin especificar"
$ws.Range("I106").Value2 = "Primera"
$ws.Range("J106").Value2 = 1300
$ws.Range("K106").Value2 = 800
$ws.Range("L106").Value2 = 900
$ws.Range("M106").Value2 = 854
$ws.Range("N106").Value2 = "`$/unidad"
$ws.Range("O106").Value2 = "Región Metropolitana"
$ws.Range("P106").Value2 = 854
$ws.Range("Q106").Value2 = 1
$ws.Range("R106").Value2 = "Hortaliza"
